# Insert 6 new weekly rows of "Choclo" price data above the existing
# row 794, pushing the previous rows 794-852 down to 800-858.
# (Commit: "Fruta / hortaliza, semanal")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at the top of the affected range.
$ws.Range("A794:A799").EntireRow.Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112024
$categoria = "Choclo"
$clasif    = "Hortaliza"

# New rows to insert (r => row number, d => date serial, h => Variedad,
# i => Calidad, j..q => Volumen..Kg o Unidades).
$newRows = @(
    @{ r = 794; d = 44578; h = "Choclero";           i = "Primera"; j = 118000; k = 150; l = 200; m = 175; n = '$/unidad'; o = "Región Metropolitana";    p = 175; q = 1 },
    @{ r = 795; d = 44578; h = "Choclero";           i = "Primera"; j = 91000;  k = 160; l = 200; m = 183; n = '$/unidad'; o = "Región de O'Higgins";     p = 183; q = 1 },
    @{ r = 796; d = 44578; h = "Choclero";           i = "Segunda"; j = 52000;  k = 100; l = 120; m = 110; n = '$/unidad'; o = "Región Metropolitana";    p = 110; q = 1 },
    @{ r = 797; d = 44578; h = "Choclero";           i = "Segunda"; j = 46000;  k = 130; l = 150; m = 143; n = '$/unidad'; o = "Región de O'Higgins";     p = 143; q = 1 },
    @{ r = 798; d = 44578; h = "Dulce o Americano";  i = "Primera"; j = 115000; k = 150; l = 170; m = 160; n = '$/unidad'; o = "Región Metropolitana";    p = 160; q = 1 },
    @{ r = 799; d = 44578; h = "Dulce o Americano";  i = "Segunda"; j = 48000;  k = 100; l = 100; m = 100; n = '$/unidad'; o = "Región Metropolitana";    p = 100; q = 1 }
)

foreach ($row in $newRows) {
    $rn = $row.r
    $ws.Cells.Item($rn, 1).Value  = $mercadoId
    $ws.Cells.Item($rn, 2).Value  = $mercado
    $ws.Cells.Item($rn, 3).Value  = $region
    $ws.Cells.Item($rn, 4).Value  = $row.d
    $ws.Cells.Item($rn, 5).Value  = $codreg
    $ws.Cells.Item($rn, 6).Value  = $catId
    $ws.Cells.Item($rn, 7).Value  = $categoria
    $ws.Cells.Item($rn, 8).Value  = $row.h
    $ws.Cells.Item($rn, 9).Value  = $row.i
    $ws.Cells.Item($rn, 10).Value = $row.j
    $ws.Cells.Item($rn, 11).Value = $row.k
    $ws.Cells.Item($rn, 12).Value = $row.l
    $ws.Cells.Item($rn, 13).Value = $row.m
    $ws.Cells.Item($rn, 14).Value = $row.n
    $ws.Cells.Item($rn, 15).Value = $row.o
    $ws.Cells.Item($rn, 16).Value = $row.p
    $ws.Cells.Item($rn, 17).Value = $row.q
    $ws.Cells.Item($rn, 18).Value = $clasif
}
